# Auto-generated script to apply cryptos list update (Sep 18 2023 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated "Price" values are purely numeric strings (e.g. "19.85"). Excel would
# normally auto-detect and store these as numbers, but the source data must stay as
# text (matching the original inline-string cells, e.g. "27.378.68" or "0.0629").
# Temporarily force a Text number format on those specific cells so the assigned
# string is kept verbatim, then restore the default "Normal" style so no stray
# formatting is left behind.
$forceTextCells = @(
    "D5", "D6", "D9", "D10", "D14", "D15", "D16", "D19", "D21", "D22",
    "D23", "D24", "D27", "D28", "D29", "D31", "D32", "D36", "D37", "D38",
    "D39", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50"
)
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.378.68"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "1.666.02"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  -0.85%  "
$ws.Range("D5").Value = "220.41"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "0.505"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("E8").Value = "  +1.24%  "
$ws.Range("D9").Value = "0.0629"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "19.85"
$ws.Range("E10").Value = "  +3.40%  "
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").Value = "1.898.56"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.650.76"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "4.21"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "0.534"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "66.95"
$ws.Range("E16").Value = "  +3.48%  "
$ws.Range("D17").Value = "27.360.32"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").Value = "224.13"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "6.77"
$ws.Range("E21").Value = "  +8.12%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.44"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "2.43"
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("D24").Value = "9.31"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "7.42"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("D28").Value = "0.120"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "16.02"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "1.21"
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("D35").Value = "1.273.12"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("D36").Value = "2.44"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "0.0177"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D38").Value = "0.540"
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("D39").Value = "0.836"
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "5.40"
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("D43").Value = "1.810.01"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "2.14"
$ws.Range("E44").Value = "  -4.48%  "
$ws.Range("D45").Value = "62.22"
$ws.Range("E45").Value = "  +1.07%  "
$ws.Range("D46").Value = "92.75"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "1.63"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").Value = "0.0517"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.70"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0983"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("E51").Value = "  +0.00%  "

# Restore default styling on the cells we temporarily reformatted
foreach ($ref in $forceTextCells) {
    $ws.Range($ref).Style = "Normal"
}

Write-Host "Applied 96 cell updates"
